# Edit script: reproduces the two substantive changes from the commit
# "Sun, May 10, 2020  3:06:24 AM":
#
#   1. The table on slide 6 (the "Sources of finance" table) is switched
#      to a different built-in PowerPoint table style - GUID
#      {F5F0D83F-1D86-455C-8CCF-22FC7855BDB4} -> {C8984D3D-F08C-4562-B7D9-AC52EBF057B2}.
#
#   2. The deck's theme (ppt/theme/theme1.xml, used by the slide master)
#      is switched from the custom "Integral" palette to the stock
#      Office default palette (name "Office" / "Office Theme"); the font
#      scheme and format scheme are already identical between the two
#      themes, so only the 12 theme colors (and the theme's display
#      name, best effort) need to change.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$slide6 = $p.Slides.Item(6)
$tableShape = $slide6.Shapes.Item(2)
if ($tableShape.HasTable) {
    $tableShape.Table.ApplyStyle("{C8984D3D-F08C-4562-B7D9-AC52EBF057B2}")
}

# --- 2. Swap the slide-master theme palette to the Office default ---------
$master = $p.SlideMaster
$theme  = $master.Theme

# Best-effort rename (some hosts don't persist this, but it's harmless).
try { $theme.Name = "Office Theme" } catch {}

$colorScheme = $theme.ThemeColorScheme
try { $colorScheme.Name = "Office" } catch {}

# Index -> (theme slot, target RGB as 0xBBGGRR for the COM .RGB property)
# matching <a:clrScheme name="Office"> from the target theme XML:
#   1 dk1      000000
#   2 lt1      FFFFFF
#   3 dk2      44546A
#   4 lt2      E7E6E6
#   5 accent1  5B9BD5
#   6 accent2  ED7D31
#   7 accent3  A5A5A5
#   8 accent4  FFC000
#   9 accent5  4472C4
#  10 accent6  70AD47
#  11 hlink    0563C1
#  12 folHlink 954F72
$officeColors = @(
    0,           # dk1      000000
    16777215,    # lt1      FFFFFF
    6968388,     # dk2      44546A
    15132391,    # lt2      E7E6E6
    13998939,    # accent1  5B9BD5
    3243501,     # accent2  ED7D31
    10855845,    # accent3  A5A5A5
    49407,       # accent4  FFC000
    12874308,    # accent5  4472C4
    4697456,     # accent6  70AD47
    12673797,    # hlink    0563C1
    7491477      # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
